# =====================================================================
# Edit script: adds "Player Info" and "ODI Batting Extra" sheets,
# renames MATCH_CARD_LINK -> MATCH_CODE (storing bare match codes
# instead of full scorecard URLs) on the existing "ODI Batting" and
# "ODI Bowling" sheets.
#
# NOTE: worksheet object references captured before a Worksheets.Add()
# call can become stale/repointed afterwards (the engine appears to
# track them positionally), so every sheet reference used for actual
# cell writes is re-fetched by name immediately before use.
# =====================================================================

function Set-TextValue($cell, [string]$text) {
    # Force the cell to be written as text (not auto-converted to a
    # number) and then drop back to the default "Normal" style so we
    # don't leave stray formatting behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-HeaderCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
}

function Get-MatchCodeFromUrl([string]$url) {
    $marker = "MatchCode="
    $idx = $url.IndexOf($marker)
    if ($idx -ge 0) {
        return $url.Substring($idx + $marker.Length)
    }
    return $url
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# First, create the new sheets in the right positions so the final
# sheet order is: Player Info, ODI Batting, ODI Bowling, ODI Batting Extra
# ---------------------------------------------------------------
$battingSheetForInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetForInsert)
$playerInfo.Name = "Player Info"

$bowlingSheetForInsert = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowlingSheetForInsert)
$extra.Name = "ODI Batting Extra"

# ---------------------------------------------------------------
# 1) Populate "Player Info"
# ---------------------------------------------------------------
$playerInfo = $wb.Worksheets.Item("Player Info")

Set-HeaderCell $playerInfo.Cells.Item(1,1) "ID"
Set-HeaderCell $playerInfo.Cells.Item(1,2) "NAME"
Set-HeaderCell $playerInfo.Cells.Item(1,3) "BATTING_HAND"
Set-HeaderCell $playerInfo.Cells.Item(1,4) "BOWL_STYLE"

Set-TextValue $playerInfo.Cells.Item(2,1) "3898"
Set-TextValue $playerInfo.Cells.Item(2,2) "Yasir Shah"
Set-TextValue $playerInfo.Cells.Item(2,3) "Right Handed"
Set-TextValue $playerInfo.Cells.Item(2,4) "Right Arm Leg Break"

# ---------------------------------------------------------------
# 2) "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
#    Replace each scorecard URL with the bare match code, and drop
#    the now-unused, already-empty INNING_NUMBER placeholder cells.
# ---------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")

Set-TextValue $battingSheet.Cells.Item(1,4) "MATCH_CODE"

$battingLastRow = 26
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $url = [string]$cell.Value2
    if ($url -ne "") {
        $code = Get-MatchCodeFromUrl($url)
        Set-TextValue $cell $code
    }
}

$emptyInningRows = @(2,4,5,7,9,10,12,14,19,20,21,26)
foreach ($r in $emptyInningRows) {
    $battingSheet.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------
# 3) "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

Set-TextValue $bowlingSheet.Cells.Item(1,2) "MATCH_CODE"

$bowlingLastRow = 25
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $url = [string]$cell.Value2
    if ($url -ne "") {
        $code = Get-MatchCodeFromUrl($url)
        Set-TextValue $cell $code
    }
}

# ---------------------------------------------------------------
# 4) Populate "ODI Batting Extra"
# ---------------------------------------------------------------
$extra = $wb.Worksheets.Item("ODI Batting Extra")

Set-HeaderCell $extra.Cells.Item(1,1) "MATCH_CODE"
Set-HeaderCell $extra.Cells.Item(1,2) "BATTING_POSITION"
Set-HeaderCell $extra.Cells.Item(1,3) "NUM_4"
Set-HeaderCell $extra.Cells.Item(1,4) "NUM_6"
Set-HeaderCell $extra.Cells.Item(1,5) "PERCENT_RUNS_OF_TOTAL"
Set-HeaderCell $extra.Cells.Item(1,6) "MAN_OF_MATCH"

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("3814", $null, $null, $null, $null, "NO"),
    @("3819", $null, $null, $null, $null, "NO"),
    @("3820", $null, $null, $null, $null, "NO"),
    @("3821", 9,     $null, $null, $null, "NO"),
    @("3822", $null, $null, $null, $null, "NO"),
    @("3836", 10,    $null, $null, $null, "YES"),
    @("3837", $null, $null, $null, $null, "NO"),
    @("3858", $null, $null, $null, $null, "NO"),
    @("3859", 10,    "1",   "1",   "8.51%", "NO"),
    @("3863", 10,    "0",   "0",   "1.85%", "NO"),
    @("3926", 10,    "0",   "0",   $null,   "NO"),
    @("3928", 10,    "3",   "0",   "9.45%", "NO"),
    @("4176", 9,     $null, $null, $null, "NO"),
    @("4177", 9,     $null, $null, $null, "NO"),
    @("4273", 10,    $null, $null, $null, "NO"),
    @("4274", 9,     "0",   "0",   "0.35%", "NO"),
    @("4275", 8,     "0",   "0",   "5.38%", "NO"),
    @("4276", $null, $null, $null, $null, "NO"),
    @("4277", 8,     "0",   "0",   "3.58%", "NO"),
    @("4292", $null, $null, $null, $null, $null)
)

$rowIdx = 2
foreach ($rowData in $extraRows) {
    Set-TextValue $extra.Cells.Item($rowIdx, 1) $rowData[0]

    if ($null -ne $rowData[1]) {
        $extra.Cells.Item($rowIdx, 2).Value = $rowData[1]
    }
    if ($null -ne $rowData[2]) {
        Set-TextValue $extra.Cells.Item($rowIdx, 3) $rowData[2]
    }
    if ($null -ne $rowData[3]) {
        Set-TextValue $extra.Cells.Item($rowIdx, 4) $rowData[3]
    }
    if ($null -ne $rowData[4]) {
        Set-TextValue $extra.Cells.Item($rowIdx, 5) $rowData[4]
    }
    if ($null -ne $rowData[5]) {
        Set-TextValue $extra.Cells.Item($rowIdx, 6) $rowData[5]
    }

    $rowIdx++
}

# ---------------------------------------------------------------
# Leave the workbook focused on the first sheet, matching activeTab=0
# ---------------------------------------------------------------
$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Activate()
$playerInfo.Range("A1").Select()
